{"js": "// Replace each three-digit-division equation in the document with its\n// updated counterpart, per the authoring diff. Every equation text in the\n// doc is unique, so an exact-match search/replace per pair is unambiguous.\nconst replacements = [\n  [\"346\u00f78=43, 2\", \"993\u00f72=496, 1\"],\n  [\"293\u00f79=32, 5\", \"296\u00f77=42, 2\"],\n  [\"347\u00f79=38, 5\", \"310\u00f78=38, 6\"],\n  [\"417\u00f74=104, 1\", \"863\u00f76=143, 5\"],\n  [\"266\u00f79=29, 5\", \"390\u00f78=48, 6\"],\n  [\"600\u00f79=66, 6\", \"805\u00f73=268, 1\"],\n  [\"157\u00f79=17, 4\", \"295\u00f76=49, 1\"],\n  [\"177\u00f72=88, 1\", \"176\u00f73=58, 2\"],\n  [\"528\u00f73=176, 0\", \"672\u00f73=224, 0\"],\n  [\"986\u00f73=328, 2\", \"427\u00f75=85, 2\"],\n  [\"884\u00f72=442, 0\", \"725\u00f76=120, 5\"],\n  [\"870\u00f75=174, 0\", \"519\u00f79=57, 6\"],\n  [\"401\u00f77=57, 2\", \"789\u00f79=87, 6\"],\n  [\"598\u00f77=85, 3\", \"737\u00f75=147, 2\"],\n  [\"758\u00f77=108, 2\", \"550\u00f77=78, 4\"],\n  [\"841\u00f75=168, 1\", \"948\u00f76=158, 0\"],\n  [\"453\u00f74=113, 1\", \"920\u00f79=102, 2\"],\n  [\"612\u00f77=87, 3\", \"930\u00f76=155, 0\"],\n  [\"446\u00f74=111, 2\", \"576\u00f72=288, 0\"],\n  [\"128\u00f72=64, 0\", \"524\u00f78=65, 4\"],\n  [\"845\u00f74=211, 1\", \"684\u00f74=171, 0\"],\n  [\"666\u00f75=133, 1\", \"110\u00f76=18, 2\"],\n  [\"701\u00f77=100, 1\", \"854\u00f74=213, 2\"],\n  [\"781\u00f73=260, 1\", \"451\u00f73=150, 1\"],\n  [\"696\u00f74=174, 0\", \"158\u00f77=22, 4\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-division equation in the document with its\n# updated counterpart, per the authoring diff. Every equation text in the\n# doc is unique, so an exact Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"346\u00f78=43, 2\", \"993\u00f72=496, 1\"),\n    @(\"293\u00f79=32, 5\", \"296\u00f77=42, 2\"),\n    @(\"347\u00f79=38, 5\", \"310\u00f78=38, 6\"),\n    @(\"417\u00f74=104, 1\", \"863\u00f76=143, 5\"),\n    @(\"266\u00f79=29, 5\", \"390\u00f78=48, 6\"),\n    @(\"600\u00f79=66, 6\", \"805\u00f73=268, 1\"),\n    @(\"157\u00f79=17, 4\", \"295\u00f76=49, 1\"),\n    @(\"177\u00f72=88, 1\", \"176\u00f73=58, 2\"),\n    @(\"528\u00f73=176, 0\", \"672\u00f73=224, 0\"),\n    @(\"986\u00f73=328, 2\", \"427\u00f75=85, 2\"),\n    @(\"884\u00f72=442, 0\", \"725\u00f76=120, 5\"),\n    @(\"870\u00f75=174, 0\", \"519\u00f79=57, 6\"),\n    @(\"401\u00f77=57, 2\", \"789\u00f79=87, 6\"),\n    @(\"598\u00f77=85, 3\", \"737\u00f75=147, 2\"),\n    @(\"758\u00f77=108, 2\", \"550\u00f77=78, 4\"),\n    @(\"841\u00f75=168, 1\", \"948\u00f76=158, 0\"),\n    @(\"453\u00f74=113, 1\", \"920\u00f79=102, 2\"),\n    @(\"612\u00f77=87, 3\", \"930\u00f76=155, 0\"),\n    @(\"446\u00f74=111, 2\", \"576\u00f72=288, 0\"),\n    @(\"128\u00f72=64, 0\", \"524\u00f78=65, 4\"),\n    @(\"845\u00f74=211, 1\", \"684\u00f74=171, 0\"),\n    @(\"666\u00f75=133, 1\", \"110\u00f76=18, 2\"),\n    @(\"701\u00f77=100, 1\", \"854\u00f74=213, 2\"),\n    @(\"781\u00f73=260, 1\", \"451\u00f73=150, 1\"),\n    @(\"696\u00f74=174, 0\", \"158\u00f77=22, 4\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
